# Daily attendance processing - 2025-12-07 16:54:29
# Normalizes the "Recorded By" column (G) so that entries beginning with a
# recorded user's email address instead lead with "System" (i.e. the first
# two comma-separated recorder names are swapped), matching the canonical
# ordering used by the attendance pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "Recorded By" (column G) whose recorder list needs the first
# two comma-separated entries swapped.
$rows = @(2,3,4,5,6,8,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,30,31,32,34,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,56,57,58,60,62,63,64,65,66,67,69,70,71,72,73,74,76,78,80,81,82,83,84,85,86,87,90,92,93,94,96,99,101,106,107,108,109,110,111,112,113,116,118,119,120,122,125,127,132,133,134,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text
    $parts = $current.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }
    if ($parts.Length -ge 2) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value = [string]::Join(", ", $parts)
    }
}
